$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 77, shifting existing rows 77-115 down to 78-116.
$ws.Rows.Item(77).Insert()

# Match the date-column number format used by the rest of the table (column D).
$ws.Cells.Item(77, 4).NumberFormat = $ws.Cells.Item(78, 4).NumberFormat

# Populate the new row 77 with its data.
$ws.Cells.Item(77, 1).Value = 8
$ws.Cells.Item(77, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(77, 3).Value = "Coquimbo"
$ws.Cells.Item(77, 4).Value = 45007
$ws.Cells.Item(77, 5).Value = 4
$ws.Cells.Item(77, 6).Value = 100114007
$ws.Cells.Item(77, 7).Value = "Jengibre"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 200
$ws.Cells.Item(77, 11).Value = 16500
$ws.Cells.Item(77, 12).Value = 17000
$ws.Cells.Item(77, 13).Value = 16750
$ws.Cells.Item(77, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(77, 15).Value = "Perú"
$ws.Cells.Item(77, 16).Value = 1288
$ws.Cells.Item(77, 17).Value = 13
$ws.Cells.Item(77, 18).Value = "Hortaliza"
